$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.792.73"
$ws.Range("E2").Value = "  +3.09%  "
$ws.Range("D3").Value = "2.627.74"
$ws.Range("E3").Value = "  +4.88%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Formula = "'326.63"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Formula = "'110.21"
$ws.Range("E6").Value = "  +1.99%  "
$ws.Range("D7").Formula = "'0.536"
$ws.Range("E7").Value = "  +1.28%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Formula = "'0.561"
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("D10").Formula = "'40.88"
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").Formula = "'20.72"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").Formula = "'7.31"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "3.040.68"
$ws.Range("E15").Value = "  +4.93%  "
$ws.Range("D16").Value = "2.616.52"
$ws.Range("E16").Value = "  +4.43%  "
$ws.Range("D17").Formula = "'0.876"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "49.773.75"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("D19").Formula = "'3.12"
$ws.Range("E19").Value = "  +12.13%  "
$ws.Range("D20").Formula = "'13.38"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Formula = "'6.84"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("D22").Value = "0.0₃0956"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Formula = "'281.53"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Formula = "'72.82"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("D25").Formula = "'2.60"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Formula = "'26.68"
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("D31").Formula = "'36.21"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("D32").Formula = "'49.83"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").Formula = "'19.82"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").Formula = "'5.46"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Formula = "'0.0795"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  +5.48%  "
$ws.Range("D38").Formula = "'4.75"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Formula = "'3.10"
$ws.Range("E39").Value = "  +6.51%  "
$ws.Range("D40").Formula = "'22.77"
$ws.Range("E40").Value = "  +5.40%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Formula = "'123.57"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Formula = "'0.113"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Formula = "'2.23"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("D45").Formula = "'3.38"
$ws.Range("E45").Value = "  +6.87%  "
$ws.Range("D46").Value = "2.057.26"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").Formula = "'2.22"
$ws.Range("E47").Value = "  +11.65%  "
$ws.Range("E48").Value = "  +9.32%  "
$ws.Range("D49").Formula = "'9.03"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").Formula = "'82.08"
$ws.Range("E51").Value = "  +1.85%  "
